$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: C34 88 -> 89, E34 recalculated (89/2256)
$ws.Range("C34").Value = 89
$ws.Range("E34").Value = 0.03945035460992908

# Row 36: C36 144 -> 145, E36 recalculated (145/1930)
$ws.Range("C36").Value = 145
$ws.Range("E36").Value = 0.07512953367875648

# Row 37: C37 948 -> 955, D37 948 -> 955 (E37 stays 1)
$ws.Range("C37").Value = 955
$ws.Range("D37").Value = 955
